$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Salary" column header
$ws.Range("D1").Value = "Salary"

# --- Column A (Test Case ID) becomes quote-prefixed text "1","2","3" ---
# Set the first cell's value + number format, then propagate that exact
# style (quote-prefix + text format) to the rest of the column via a
# format-only paste so no stray intermediate styles are produced.
$ws.Range("A2").Value = "'1"
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3:A4").PasteSpecial(-4122) | Out-Null
$ws.Range("A3").Value = "'2"
$ws.Range("A4").Value = "'3"

# --- New column D (Salary) as quote-prefixed text "5000","10000","23000" ---
$ws.Range("D2:D4").PasteSpecial(-4122) | Out-Null
$ws.Range("D2").Value = "'5000"
$ws.Range("D3").Value = "'10000"
$ws.Range("D4").Value = "'23000"

$excel.CutCopyMode = 0

$ws.Range("D5").Select()
